$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the image reference used for the "Coloana de vopsit in 3 culori" row
$ws.Range("B11").Value = "assets/image9.jpg"

# Fix typo in the NCS color code for "Negativul de la Hol" row (B08G -> B80G)
$ws.Range("E13").Value = "NCS S 2050- B80G /  NCS S 3040-B20G"

# Update the view's scroll position / active selection to match the saved state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I11").Select()
